# Updates data cells in several sheets to reflect refreshed market-price
# figures pulled by the scheduled runner (currentAveragePrice* / LevePrice* /
# LeveProfit* columns, H:N) for the Bahamut_Profits workbook.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value2 = 99.5
$ws.Range("I12").Value2 = 0
$ws.Range("J12").Value2 = 99.5
$ws.Range("K12").Value2 = 0
$ws.Range("L12").Value2 = 99.5
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value2 = -439.5
$ws.Range("H41").Value2 = 770.2222
$ws.Range("I41").Value2 = 187.75
$ws.Range("J41").Value2 = 1236.2
$ws.Range("K41").Value2 = 187.75
$ws.Range("L41").Value2 = 1236.2
$ws.Range("M41").Value2 = 252.25
$ws.Range("N41").Value2 = -2116.2
$ws.Range("H52").Value2 = 0
$ws.Range("J52").Value2 = 0
$ws.Range("L52").Value2 = 0
$ws.Range("N52").ClearContents()
$ws.Range("H88").Value2 = 1135599.8
$ws.Range("I88").Value2 = 2499
$ws.Range("J88").Value2 = 2268700.5
$ws.Range("K88").Value2 = 2499
$ws.Range("L88").Value2 = 2268700.5
$ws.Range("M88").Value2 = -2093
$ws.Range("N88").Value2 = -2269512.5
$ws.Range("H91").Value2 = 1135599.8
$ws.Range("I91").Value2 = 2499
$ws.Range("J91").Value2 = 2268700.5
$ws.Range("K91").Value2 = 2499
$ws.Range("L91").Value2 = 2268700.5
$ws.Range("M91").Value2 = -1095
$ws.Range("N91").Value2 = -2271508.5
$ws.Range("H111").Value2 = 1825.8
$ws.Range("I111").Value2 = 1782.25
$ws.Range("J111").Value2 = 2000
$ws.Range("K111").Value2 = 5346.75
$ws.Range("L111").Value2 = 6000
$ws.Range("M111").Value2 = -2279.75
$ws.Range("N111").Value2 = -12134
$ws.Range("H132").Value2 = 4582.273
$ws.Range("I132").Value2 = 7234.1665
$ws.Range("J132").Value2 = 1400
$ws.Range("K132").Value2 = 21702.4995
$ws.Range("L132").Value2 = 4200
$ws.Range("M132").Value2 = -19172.4995
$ws.Range("N132").Value2 = -9260
$ws.Range("H137").Value2 = 1118.3636
$ws.Range("I137").Value2 = 1110.2
$ws.Range("K137").Value2 = 3330.6
$ws.Range("M137").Value2 = -780.6000000000004

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 6998.8086
$ws.Range("I32").Value2 = 4729.5454
$ws.Range("J32").Value2 = 16599.54
$ws.Range("K32").Value2 = 4729.5454
$ws.Range("L32").Value2 = 16599.54
$ws.Range("M32").Value2 = -4442.5454
$ws.Range("N32").Value2 = -17173.54

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 2181.018
$ws.Range("I31").Value2 = 2583.8484
$ws.Range("J31").Value2 = 1576.7727
$ws.Range("K31").Value2 = 2583.8484
$ws.Range("L31").Value2 = 1576.7727
$ws.Range("M31").Value2 = -2288.8484
$ws.Range("N31").Value2 = -2166.7727
$ws.Range("H34").Value2 = 2181.018
$ws.Range("I34").Value2 = 2583.8484
$ws.Range("J34").Value2 = 1576.7727
$ws.Range("K34").Value2 = 2583.8484
$ws.Range("L34").Value2 = 1576.7727
$ws.Range("M34").Value2 = -2381.8484
$ws.Range("N34").Value2 = -1980.7727
$ws.Range("H86").Value2 = 12188.5
$ws.Range("I86").Value2 = 12179.4
$ws.Range("J86").Value2 = 12197.6
$ws.Range("K86").Value2 = 12179.4
$ws.Range("L86").Value2 = 12197.6
$ws.Range("M86").Value2 = -11056.4
$ws.Range("N86").Value2 = -14443.6
$ws.Range("H89").Value2 = 12188.5
$ws.Range("I89").Value2 = 12179.4
$ws.Range("J89").Value2 = 12197.6
$ws.Range("K89").Value2 = 60897
$ws.Range("L89").Value2 = 60988
$ws.Range("M89").Value2 = -55281
$ws.Range("N89").Value2 = -72220
$ws.Range("H134").Value2 = 31251992
$ws.Range("I134").Value2 = 2104.7856
$ws.Range("J134").Value2 = 250001200
$ws.Range("K134").Value2 = 6314.3568
$ws.Range("L134").Value2 = 750003600
$ws.Range("M134").Value2 = -3779.3568
$ws.Range("N134").Value2 = -750008670

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 3460405
$ws.Range("I2").Value2 = 8403532
$ws.Range("J2").Value2 = 216
$ws.Range("K2").Value2 = 50421192
$ws.Range("L2").Value2 = 1296
$ws.Range("M2").Value2 = -50421079
$ws.Range("N2").Value2 = -1522
$ws.Range("H5").Value2 = 1218.9412
$ws.Range("I5").Value2 = 1567.6364
$ws.Range("K5").Value2 = 4702.9092
$ws.Range("M5").Value2 = -4590.9092
$ws.Range("H13").Value2 = 175
$ws.Range("I13").Value2 = 0
$ws.Range("K13").Value2 = 0
$ws.Range("M13").ClearContents()
$ws.Range("H41").Value2 = 10003
$ws.Range("J41").Value2 = 10003
$ws.Range("L41").Value2 = 30009
$ws.Range("N41").Value2 = -30685
$ws.Range("H112").Value2 = 2486.7273
$ws.Range("I112").Value2 = 1336.2858
$ws.Range("K112").Value2 = 4008.8574
$ws.Range("M112").Value2 = -2900.8574
$ws.Range("H113").Value2 = 588.4286
$ws.Range("I113").Value2 = 990
$ws.Range("J113").Value2 = 550.78125
$ws.Range("K113").Value2 = 2970
$ws.Range("L113").Value2 = 1652.34375
$ws.Range("M113").Value2 = -800
$ws.Range("N113").Value2 = -5992.34375
$ws.Range("H131").Value2 = 6593115
$ws.Range("J131").Value2 = 952.7941
$ws.Range("L131").Value2 = 2858.3823
$ws.Range("N131").Value2 = -12938.3823
$ws.Range("H135").Value2 = 1218.9412
$ws.Range("I135").Value2 = 1567.6364
$ws.Range("K135").Value2 = 14108.7276
$ws.Range("M135").Value2 = -11573.7276

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value2 = 64.8
$ws.Range("I2").Value2 = 51.083332
$ws.Range("J2").Value2 = 119.666664
$ws.Range("K2").Value2 = 51.083332
$ws.Range("L2").Value2 = 119.666664
$ws.Range("M2").Value2 = 61.916668
$ws.Range("N2").Value2 = -345.666664

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value2 = 5000
$ws.Range("J121").Value2 = 5000
$ws.Range("L121").Value2 = 5000
$ws.Range("N121").Value2 = -8494
$ws.Range("H122").Value2 = 7014.1904
$ws.Range("I122").Value2 = 8787.4
$ws.Range("K122").Value2 = 26362.2
$ws.Range("M122").Value2 = -23912.2
$ws.Range("H132").Value2 = 3561
$ws.Range("I132").Value2 = 2622.5
$ws.Range("J132").Value2 = 4499.5
$ws.Range("K132").Value2 = 7867.5
$ws.Range("L132").Value2 = 13498.5
$ws.Range("M132").Value2 = -5337.5
$ws.Range("N132").Value2 = -18558.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value2 = 38643
$ws.Range("J110").Value2 = 38643
$ws.Range("L110").Value2 = 38643
$ws.Range("N110").Value2 = -46823
$ws.Range("H112").Value2 = 25500
$ws.Range("J112").Value2 = 25500
$ws.Range("L112").Value2 = 25500
$ws.Range("N112").Value2 = -28454
